$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New positional ranking values for B2:K12 (columns 2-11, rows 2-12)
$data = @(
  @(2, 0, 0, 0, 3, 1, 0, 0, 0, 0),
  @(1, 0, 0, 2, 0, 4, 0, 0, 0, 0),
  @(0, 0, 1, 0, 0, 0, 4, 0, 0, 0),
  @(0, 0, 0, 2, 0, 0, 5, 5, 0, 0),
  @(0, 0, 10, 10, 2, 3, 2, 2, 2, 2),
  @(0, 0, 10, 0, 1, 0, 0, 0, 0, 0),
  @(0, 0, 0, 4, 0, 0, 1, 0, 0, 0),
  @(0, 0, 10, 10, 0, 1, 5, 3, 1, 2),
  @(0, 0, 10, 4, 3, 3, 10, 1, 2, 2),
  @(0, 0, 0, 0, 0, 0, 10, 5, 4, 1),
  @(4, 4, 4, 4, 4, 4, 4, 4, 4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = 2 + $i
  $rowValues = $data[$i]
  for ($j = 0; $j -lt $rowValues.Length; $j++) {
    $col = 2 + $j
    $ws.Cells.Item($row, $col).Value = $rowValues[$j]
  }
}

# Select cell H15 to match the saved selection state
$ws.Range("H15").Select()
